$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of tracker data
$ws.Range("A3").Value = "G2"
$ws.Range("B3").Value = "sedrftgyhuioygtfrd"

$ws.Range("C3").Value = 45889
$ws.Range("C3").Style = $ws.Range("C2").Style
$ws.Range("C3").NumberFormat = $ws.Range("C2").NumberFormat

$ws.Range("D3").Value = 1.01
$ws.Range("E3").Value = 100
$ws.Range("F3").Value = 0.01
